# Auto-generated edit script: updates market-price-derived columns (H-N)
# on specific Leve rows across multiple job sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 7353.8096
$ws.Cells.Item(53, 9).Value = 11637.538
$ws.Cells.Item(53, 10).Value = 392.75
$ws.Cells.Item(53, 11).Value = 11637.538
$ws.Cells.Item(53, 12).Value = 392.75
$ws.Cells.Item(53, 13).Value = -11000.538
$ws.Cells.Item(53, 14).Value = -1666.75

$ws.Cells.Item(69, 8).Value = 6528.2354
$ws.Cells.Item(69, 9).Value = 9000
$ws.Cells.Item(69, 10).Value = 4798
$ws.Cells.Item(69, 11).Value = 27000
$ws.Cells.Item(69, 12).Value = 14394
$ws.Cells.Item(69, 13).Value = -26126
$ws.Cells.Item(69, 14).Value = -16142

$ws.Cells.Item(72, 8).Value = 6528.2354
$ws.Cells.Item(72, 9).Value = 9000
$ws.Cells.Item(72, 10).Value = 4798
$ws.Cells.Item(72, 11).Value = 81000
$ws.Cells.Item(72, 12).Value = 43182
$ws.Cells.Item(72, 13).Value = -76632
$ws.Cells.Item(72, 14).Value = -51918

$ws.Cells.Item(76, 8).Value = 43481812
$ws.Cells.Item(76, 9).Value = 58827204
$ws.Cells.Item(76, 10).Value = 3206.6667
$ws.Cells.Item(76, 11).Value = 58827204
$ws.Cells.Item(76, 12).Value = 3206.6667
$ws.Cells.Item(76, 13).Value = -58826889
$ws.Cells.Item(76, 14).Value = -3836.6667

$ws.Cells.Item(79, 8).Value = 43481812
$ws.Cells.Item(79, 9).Value = 58827204
$ws.Cells.Item(79, 10).Value = 3206.6667
$ws.Cells.Item(79, 11).Value = 58827204
$ws.Cells.Item(79, 12).Value = 3206.6667
$ws.Cells.Item(79, 13).Value = -58826112
$ws.Cells.Item(79, 14).Value = -5390.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1861.9231
$ws.Cells.Item(63, 9).Value = 1814.762
$ws.Cells.Item(63, 10).Value = 2060
$ws.Cells.Item(63, 11).Value = 1814.762
$ws.Cells.Item(63, 12).Value = 2060
$ws.Cells.Item(63, 13).Value = -1128.762
$ws.Cells.Item(63, 14).Value = -3432

$ws.Cells.Item(66, 8).Value = 1861.9231
$ws.Cells.Item(66, 9).Value = 1814.762
$ws.Cells.Item(66, 10).Value = 2060
$ws.Cells.Item(66, 11).Value = 9073.809999999999
$ws.Cells.Item(66, 12).Value = 10300
$ws.Cells.Item(66, 13).Value = -5641.809999999999
$ws.Cells.Item(66, 14).Value = -17164

$ws.Cells.Item(74, 8).Value = 33334858
$ws.Cells.Item(74, 9).Value = 31250698
$ws.Cells.Item(74, 11).Value = 31250698
$ws.Cells.Item(74, 13).Value = -31249824

$ws.Cells.Item(77, 8).Value = 33334858
$ws.Cells.Item(77, 9).Value = 31250698
$ws.Cells.Item(77, 11).Value = 156253490
$ws.Cells.Item(77, 13).Value = -156249122

$ws.Cells.Item(132, 8).Value = 19590126
$ws.Cells.Item(132, 9).Value = 21289142
$ws.Cells.Item(132, 11).Value = 63867426
$ws.Cells.Item(132, 13).Value = -63864896

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2140658
$ws.Cells.Item(31, 9).Value = 2977898.2
$ws.Cells.Item(31, 10).Value = 9501.182000000001
$ws.Cells.Item(31, 11).Value = 2977898.2
$ws.Cells.Item(31, 12).Value = 9501.182000000001
$ws.Cells.Item(31, 13).Value = -2977603.2
$ws.Cells.Item(31, 14).Value = -10091.182

$ws.Cells.Item(34, 8).Value = 2140658
$ws.Cells.Item(34, 9).Value = 2977898.2
$ws.Cells.Item(34, 10).Value = 9501.182000000001
$ws.Cells.Item(34, 11).Value = 2977898.2
$ws.Cells.Item(34, 12).Value = 9501.182000000001
$ws.Cells.Item(34, 13).Value = -2977696.2
$ws.Cells.Item(34, 14).Value = -9905.182000000001

$ws.Cells.Item(70, 8).Value = 21600
$ws.Cells.Item(70, 10).Value = 21600
$ws.Cells.Item(70, 12).Value = 21600
$ws.Cells.Item(70, 14).Value = -22230

$ws.Cells.Item(73, 8).Value = 21600
$ws.Cells.Item(73, 10).Value = 21600
$ws.Cells.Item(73, 12).Value = 21600
$ws.Cells.Item(73, 14).Value = -23784

$ws.Cells.Item(132, 8).Value = 2153.2258
$ws.Cells.Item(132, 9).Value = 1455.75
$ws.Cells.Item(132, 11).Value = 4367.25
$ws.Cells.Item(132, 13).Value = -1837.25

$ws.Cells.Item(134, 8).Value = 1030309.5
$ws.Cells.Item(134, 9).Value = 5823
$ws.Cells.Item(134, 10).Value = 2503008.8
$ws.Cells.Item(134, 11).Value = 17469
$ws.Cells.Item(134, 12).Value = 7509026.399999999
$ws.Cells.Item(134, 13).Value = -14934
$ws.Cells.Item(134, 14).Value = -7514096.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1454665.6
$ws.Cells.Item(5, 9).Value = 560.5417
$ws.Cells.Item(5, 10).Value = 3291430
$ws.Cells.Item(5, 11).Value = 1681.6251
$ws.Cells.Item(5, 12).Value = 9874290
$ws.Cells.Item(5, 13).Value = -1569.6251
$ws.Cells.Item(5, 14).Value = -9874514

$ws.Cells.Item(106, 8).Value = 4935.2383
$ws.Cells.Item(106, 10).Value = 4935.2383
$ws.Cells.Item(106, 12).Value = 14805.7149
$ws.Cells.Item(106, 14).Value = -16697.7149

$ws.Cells.Item(132, 8).Value = 1315.6177
$ws.Cells.Item(132, 9).Value = 694.4375
$ws.Cells.Item(132, 11).Value = 6249.9375
$ws.Cells.Item(132, 13).Value = -3719.9375

$ws.Cells.Item(135, 8).Value = 1454665.6
$ws.Cells.Item(135, 9).Value = 560.5417
$ws.Cells.Item(135, 10).Value = 3291430
$ws.Cells.Item(135, 11).Value = 5044.8753
$ws.Cells.Item(135, 12).Value = 29622870
$ws.Cells.Item(135, 13).Value = -2509.8753
$ws.Cells.Item(135, 14).Value = -29627940

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6043966
$ws.Cells.Item(70, 9).Value = 2608396.5
$ws.Cells.Item(70, 10).Value = 14289333
$ws.Cells.Item(70, 11).Value = 2608396.5
$ws.Cells.Item(70, 12).Value = 14289333
$ws.Cells.Item(70, 13).Value = -2608126.5
$ws.Cells.Item(70, 14).Value = -14289873

$ws.Cells.Item(73, 8).Value = 6043966
$ws.Cells.Item(73, 9).Value = 2608396.5
$ws.Cells.Item(73, 10).Value = 14289333
$ws.Cells.Item(73, 11).Value = 2608396.5
$ws.Cells.Item(73, 12).Value = 14289333
$ws.Cells.Item(73, 13).Value = -2607460.5
$ws.Cells.Item(73, 14).Value = -14291205

$ws.Cells.Item(88, 8).Value = 30000
$ws.Cells.Item(88, 10).Value = 30000
$ws.Cells.Item(88, 12).Value = 30000
$ws.Cells.Item(88, 14).Value = -30902

$ws.Cells.Item(91, 8).Value = 30000
$ws.Cells.Item(91, 10).Value = 30000
$ws.Cells.Item(91, 12).Value = 30000
$ws.Cells.Item(91, 14).Value = -33120

$ws.Cells.Item(102, 8).Value = 4425.4287
$ws.Cells.Item(102, 9).Value = 4882.4443
$ws.Cells.Item(102, 10).Value = 1683.3334
$ws.Cells.Item(102, 11).Value = 4882.4443
$ws.Cells.Item(102, 12).Value = 1683.3334
$ws.Cells.Item(102, 13).Value = -3260.4443
$ws.Cells.Item(102, 14).Value = -4927.3334

$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2209.6296
$ws.Cells.Item(61, 9).Value = 1603.75
$ws.Cells.Item(61, 10).Value = 3090.9092
$ws.Cells.Item(61, 11).Value = 1603.75
$ws.Cells.Item(61, 12).Value = 3090.9092
$ws.Cells.Item(61, 13).Value = -1401.75
$ws.Cells.Item(61, 14).Value = -3494.9092

$ws.Cells.Item(113, 8).Value = 2209.6296
$ws.Cells.Item(113, 9).Value = 1603.75
$ws.Cells.Item(113, 10).Value = 3090.9092
$ws.Cells.Item(113, 11).Value = 1603.75
$ws.Cells.Item(113, 12).Value = 3090.9092
$ws.Cells.Item(113, 13).Value = 566.25
$ws.Cells.Item(113, 14).Value = -7430.9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 21885
$ws.Cells.Item(70, 9).Value = 25000
$ws.Cells.Item(70, 10).Value = 21365.834
$ws.Cells.Item(70, 11).Value = 25000
$ws.Cells.Item(70, 12).Value = 21365.834
$ws.Cells.Item(70, 13).Value = -24685
$ws.Cells.Item(70, 14).Value = -21995.834

$ws.Cells.Item(73, 8).Value = 21885
$ws.Cells.Item(73, 9).Value = 25000
$ws.Cells.Item(73, 10).Value = 21365.834
$ws.Cells.Item(73, 11).Value = 25000
$ws.Cells.Item(73, 12).Value = 21365.834
$ws.Cells.Item(73, 13).Value = -23908
$ws.Cells.Item(73, 14).Value = -23549.834

$ws.Cells.Item(122, 8).Value = 1950533.9
$ws.Cells.Item(122, 9).Value = 2646429.2
$ws.Cells.Item(122, 10).Value = 2026.5333
$ws.Cells.Item(122, 11).Value = 7939287.600000001
$ws.Cells.Item(122, 12).Value = 6079.5999
$ws.Cells.Item(122, 13).Value = -7936837.600000001

$ws.Cells.Item(136, 8).Value = 1690.6875
$ws.Cells.Item(136, 9).Value = 1293.5652
$ws.Cells.Item(136, 10).Value = 1913.4634
$ws.Cells.Item(136, 11).Value = 3880.6956
$ws.Cells.Item(136, 12).Value = 5740.3902
$ws.Cells.Item(136, 13).Value = -1330.6956
$ws.Cells.Item(136, 14).Value = -10840.3902
